# Generate Report for Handoff
# Updates the Priority and Latest Handoff Datetime for the
# 3168d6f9-27e9-4f30-9726-229db29279e2 (.md) file across the
# zh-cn and de-de localization sheets, and refreshes the
# "Latest HO Xliff Generate Date" on the Overview sheet for that row.

$wb = $excel.ActiveWorkbook

# --- Overview sheet: row 4 is the 3168d6f9-... file ---
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("G4").Value = "2016-08-15 16:30:35"

# --- zh-cn sheet: rows 4-7 correspond to the files that were "low" priority ---
$wsZhCn = $wb.Worksheets.Item("zh-cn")
foreach ($r in 4..7) {
    $wsZhCn.Cells.Item($r, 5).Value = "ht"
    $wsZhCn.Cells.Item($r, 8).Value = "2016-08-15 16:30:31"
}

# --- de-de sheet: rows 4-7 correspond to the files that were "low" priority ---
$wsDeDe = $wb.Worksheets.Item("de-de")
foreach ($r in 4..7) {
    $wsDeDe.Cells.Item($r, 5).Value = "ht"
    $wsDeDe.Cells.Item($r, 8).Value = "2016-08-15 16:30:35"
}
